$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.652.28"
$ws.Range("E2").Value = "  -2.98%  "

$ws.Range("D3").Value = "2.899.91"
$ws.Range("E3").Value = "  -3.98%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "585.36"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").Value = "146.92"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  -3.04%  "

$ws.Range("D9").Value = "2.898.02"
$ws.Range("E9").Value = "  -3.97%  "

$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  +4.45%  "

$ws.Range("E11").Value = "  -4.40%  "

$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("D14").Value = "33.88"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").Value = "3.376.69"
$ws.Range("E16").Value = "  -4.09%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.80"
$ws.Range("E17").Value = "  -2.90%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "60.557.60"
$ws.Range("E18").Value = "  -3.10%  "

$ws.Range("D19").Value = "2.894.84"
$ws.Range("E19").Value = "  -4.36%  "

$ws.Range("D20").Value = "425.38"
$ws.Range("E20").Value = "  -5.09%  "

$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value = "  -3.93%  "

$ws.Range("E22").Value = "  -2.69%  "

$ws.Range("D23").Value = "7.09"
$ws.Range("E23").Value = "  -4.85%  "

$ws.Range("D24").Value = "80.38"
$ws.Range("E24").Value = "  -2.43%  "

$ws.Range("D25").Value = "11.06"
$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("E26").Value = "  -0.42%  "

$ws.Range("D27").Value = "11.87"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +3.20%  "

$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E32").Value = "  -3.74%  "

$ws.Range("D33").Value = "26.46"
$ws.Range("E33").Value = "  -3.89%  "

$ws.Range("E34").Value = "  -2.57%  "

$ws.Range("D35").Value = "0.0₃0833"
$ws.Range("E35").Value = "  -2.28%  "

$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -1.95%  "

$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("D38").Value = "49.19"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("E40").Value = "  -3.49%  "

$ws.Range("E41").Value = "  -0.95%  "

$ws.Range("E42").Value = "  -3.36%  "

$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  +2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.80"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("D45").Value = "0.0347"
$ws.Range("E45").Value = "  -2.12%  "

$ws.Range("D46").Value = "370.63"

$ws.Range("D47").Value = "133.26"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").Value = "2.652.42"
$ws.Range("E48").Value = "  -3.10%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "25.48"
$ws.Range("E49").Value = "  +7.18%  "

$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("E51").Value = "  -1.20%  "
